# Refresh the cryptos worksheet with the latest scraped price/volume data.
# A handful of Price cells are plain decimal-looking strings (e.g. "1.00",
# "0.578") stored as text in the source data. Column D has no explicit
# NumberFormat (it is "General"), so assigning those through .Value would let
# Excel auto-coerce them to numbers and silently drop significant trailing
# zeros/precision (e.g. "1.00" -> 1, "2.03" -> 2.0299999999999998). Set-TextValue
# works around this the way a human would in the Excel UI: prefix with an
# apostrophe to force text entry, then ClearFormats() to drop the "quote
# prefix" cell style Excel stamps on afterwards, so no stray style survives.
function Set-TextValue($range, $text) {
    $range.Value = "'" + $text
    $range.ClearFormats()
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '89.497.33'
$ws.Range("E2").Value = '  +3.05%  '
$ws.Range("D3").Value = '3.281.59'
$ws.Range("E3").Value = '  -1.53%  '
$ws.Range("E4").Value = '  +0.14%  '
Set-TextValue $ws.Range("D5") '212.87'
$ws.Range("E5").Value = '  -3.07%  '
Set-TextValue $ws.Range("D6") '627.16'
$ws.Range("E6").Value = '  -1.95%  '
Set-TextValue $ws.Range("D7") '0.378'
$ws.Range("E7").Value = '  +16.84%  '
Set-TextValue $ws.Range("D8") '0.734'
$ws.Range("E8").Value = '  +18.04%  '
$ws.Range("E9").Value = '  +0.03%  '
$ws.Range("D10").Value = '3.278.16'
$ws.Range("E10").Value = '  -1.96%  '
Set-TextValue $ws.Range("D11") '0.578'
$ws.Range("E11").Value = '  -4.46%  '
$ws.Range("E12").Value = '  +11.10%  '
$ws.Range("E13").Value = '  -4.44%  '
Set-TextValue $ws.Range("D14") '34.18'
$ws.Range("E14").Value = '  -0.69%  '
$ws.Range("D15").Value = '3.883.09'
$ws.Range("E15").Value = '  -1.23%  '
Set-TextValue $ws.Range("D16") '5.44'
$ws.Range("E16").Value = '  +0.73%  '
$ws.Range("D17").Value = '89.191.28'
$ws.Range("E17").Value = '  +2.61%  '
$ws.Range("D18").Value = '3.296.86'
$ws.Range("E18").Value = '  -0.62%  '
Set-TextValue $ws.Range("D19") '14.10'
$ws.Range("E19").Value = '  -4.03%  '
Set-TextValue $ws.Range("D20") '3.07'
$ws.Range("E20").Value = '  -4.04%  '
Set-TextValue $ws.Range("D21") '437.52'
$ws.Range("E21").Value = '  -2.08%  '
Set-TextValue $ws.Range("D22") '8.91'
$ws.Range("E22").Value = '  -2.97%  '
Set-TextValue $ws.Range("D23") '5.35'
$ws.Range("E23").Value = '  +1.75%  '
$ws.Range("E24").Value = '  -0.44%  '
$ws.Range("E25").Value = '  -2.44%  '
Set-TextValue $ws.Range("D26") '12.23'
$ws.Range("E26").Value = '  +0.61%  '
$ws.Range("D27").Value = '3.464.60'
$ws.Range("E27").Value = '  +0.81%  '
Set-TextValue $ws.Range("D28") '76.92'
$ws.Range("E28").Value = '  -2.14%  '
$ws.Range("E29").Value = '  +2.94%  '
$ws.Range("E30").Value = '  -0.05%  '
Set-TextValue $ws.Range("D31") '0.181'
$ws.Range("E31").Value = '  +0.34%  '
Set-TextValue $ws.Range("D32") '1.00'
$ws.Range("E32").Value = '  +0.00%  '
Set-TextValue $ws.Range("D33") '8.86'
$ws.Range("E33").Value = '  -4.52%  '
Set-TextValue $ws.Range("D34") '562.76'
$ws.Range("E34").Value = '  -7.26%  '
Set-TextValue $ws.Range("D35") '1.37'
$ws.Range("E35").Value = '  -11.72%  '
$ws.Range("E36").Value = '  -4.54%  '
$ws.Range("E37").Value = '  +9.06%  '
Set-TextValue $ws.Range("D38") '0.140'
$ws.Range("E38").Value = '  -7.24%  '
Set-TextValue $ws.Range("D39") '22.71'
$ws.Range("E39").Value = '  -3.08%  '
Set-TextValue $ws.Range("D40") '21.83'
$ws.Range("E40").Value = '  +2.41%  '
$ws.Range("E41").Value = '  +0.12%  '
Set-TextValue $ws.Range("D42") '3.09'
$ws.Range("E42").Value = '  -2.86%  '
$ws.Range("B43").Value = 'PolygonEcosystemToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
Set-TextValue $ws.Range("D43") '0.401'
$ws.Range("E43").Value = '  -4.38%  '
$ws.Range("B44").Value = 'Stacks'
$ws.Range("C44").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue $ws.Range("D44") '2.03'
$ws.Range("E44").Value = '  -1.47%  '
$ws.Range("E45").Value = '  +0.04%  '
Set-TextValue $ws.Range("D46") '155.13'
$ws.Range("E46").Value = '  -0.98%  '
$ws.Range("B47").Value = 'Stellar'
$ws.Range("C47").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue $ws.Range("D47") '0.136'
$ws.Range("E47").Value = '  +20.89%  '
$ws.Range("B48").Value = 'Aave'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue $ws.Range("D48") '180.69'
$ws.Range("E48").Value = '  -4.82%  '
$ws.Range("B49").Value = 'OKB'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue $ws.Range("D49") '44.99'
$ws.Range("E49").Value = '  -1.55%  '
$ws.Range("E50").Value = '  -4.26%  '
Set-TextValue $ws.Range("D51") '4.24'
$ws.Range("E51").Value = '  -0.88%  '
